$d = $word.ActiveDocument

$d.Content.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s699817", 2)

$d.Content.Find.Execute("Ref-DJ79X2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s878264", 2)

$d.Content.Find.Execute("Ref-D4E5F6", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-f472192", 2)

$d.Content.Find.Execute("(Ref-DJ49F2)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Al-Sayed, 1998)", 2)
